$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 108 (old rows 108-126 shift down to 110-128).
$ws.Range("A108:A109").EntireRow.Insert()

# New row 108: Feria Lagunitas de Puerto Montt, Mandarina / Murcott / Primera, week of 2021-11-05.
$ws.Range("A108").Value = 4
$ws.Range("B108").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C108").Value = "Los Lagos"
$ws.Range("D108").Value = 44505
$ws.Range("E108").Value = 10
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100102
$ws.Range("H108").Value = "Cítricos"
$ws.Range("I108").Value = 100102004
$ws.Range("J108").Value = "Mandarina"
$ws.Range("K108").Value = "Murcott"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 600
$ws.Range("N108").Value = 5500
$ws.Range("O108").Value = 6000
$ws.Range("P108").Value = 5750
$ws.Range("Q108").Value = "$/caja 10 kilos"
$ws.Range("R108").Value = "Provincia de Limarí"
$ws.Range("S108").Value = 575
$ws.Range("T108").Value = 10

# New row 109: same market/date, Mandarina / Murcott / Segunda.
$ws.Range("A109").Value = 4
$ws.Range("B109").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C109").Value = "Los Lagos"
$ws.Range("D109").Value = 44505
$ws.Range("E109").Value = 10
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100102
$ws.Range("H109").Value = "Cítricos"
$ws.Range("I109").Value = 100102004
$ws.Range("J109").Value = "Mandarina"
$ws.Range("K109").Value = "Murcott"
$ws.Range("L109").Value = "Segunda"
$ws.Range("M109").Value = 300
$ws.Range("N109").Value = 4000
$ws.Range("O109").Value = 4000
$ws.Range("P109").Value = 4000
$ws.Range("Q109").Value = "$/caja 10 kilos"
$ws.Range("R109").Value = "Provincia de Limarí"
$ws.Range("S109").Value = 400
$ws.Range("T109").Value = 10
